# feat: add 2022-Q1 data
#
# Before:  2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
# After:   2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#
# The existing "总计" sheet (position 5) is repurposed to hold the new
# "2022-Q1" per-fund holdings table (its old roll-up content is preserved by
# moving it onto a freshly appended "总计" sheet with one new row inserted
# on top for the 2022-Q1 summary line).

$wb = $excel.ActiveWorkbook

# Use the "2021-Q4" sheet as a style template: its header row (B1:H1) and its
# "A" column (index cell, A2) carry the workbook's shared bold/centered/bordered
# style used throughout every sheet.
$template = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# Step 1: turn the current "总计" sheet into the new "2022-Q1" fund table.
# ---------------------------------------------------------------------
$quarterSheet = $wb.Worksheets.Item(5)
$quarterSheet.Cells.Clear()
$quarterSheet.Name = "2022-Q1"

# Pre-format the text body columns (B:G) as text so codes like "001113" and
# decimal strings like "20.79" are stored verbatim instead of being coerced
# into numbers (which would drop leading zeros / trailing zeros).
$bodyTextRange = $quarterSheet.Range("B2:G7")
$bodyTextRange.NumberFormat = "@"

# Header row: copy the shared header style, then set the Chinese captions.
$srcHeader = $template.Range("B1:H1")
$srcHeader.Copy()
$quarterSheet.Range("B1:H1").PasteSpecial(-4122)

$quarterSheet.Cells.Item(1,2).Value = "基金代码"
$quarterSheet.Cells.Item(1,3).Value = "基金名称"
$quarterSheet.Cells.Item(1,4).Value = "基金规模"
$quarterSheet.Cells.Item(1,5).Value = "股票总仓位"
$quarterSheet.Cells.Item(1,6).Value = "仓位占比"
$quarterSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$quarterSheet.Cells.Item(1,8).Value = "仓位排名"

# Body rows: A = running index (number), B-G = text fields, H = rank (number)
$fundRows = @(
    @(0, "001113", "南方大数据100指数A",           "20.79", "94.23", "2.33", "0.4844", 3),
    @(1, "014125", "华夏中证1000指数增强A",         "7.03",  "89.75", "0.85", "0.0598", 3),
    @(2, "014126", "华夏中证1000指数增强C",         "6.09",  "89.75", "0.85", "0.0518", 3),
    @(3, "009658", "汇丰晋信中小盘低波动策略股票A", "0.98",  "86.56", "1.11", "0.0109", 6),
    @(4, "004344", "南方大数据100指数C",           "0.17",  "94.23", "2.33", "0.0040", 3),
    @(5, "009775", "汇丰晋信中小盘低波动策略股票C", "0.04",  "86.56", "1.11", "0.0004", 6)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = 2 + $i
    $fundRow = $fundRows[$i]
    $quarterSheet.Cells.Item($r, 1).Value = $fundRow[0]
    for ($col = 2; $col -le 7; $col++) {
        $quarterSheet.Cells.Item($r, $col).Value = [string]$fundRow[$col - 1]
    }
    $quarterSheet.Cells.Item($r, 8).Value = $fundRow[7]
}

# Drop the temporary text-number-format now that the strings are committed.
$bodyTextRange.ClearFormats()

# Apply the shared "A" column style (bold/centered/bordered) to A2:A7.
$srcIndexCell = $template.Range("A2")
$srcIndexCell.Copy()
$quarterSheet.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet after "2022-Q1" holding the old
# roll-up table, with a new first row for the 2022-Q1 summary.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add($null, $quarterSheet)
$totalSheet.Name = "总计"

$srcTotalHeader = $template.Range("B1:D1")
$srcTotalHeader.Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)

$totalSheet.Cells.Item(1,2).Value = "日期"
$totalSheet.Cells.Item(1,3).Value = "持有数量(只)"
$totalSheet.Cells.Item(1,4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 6,  0.61),
    @(1, "2021-Q4", 1,  1.14),
    @(2, "2021-Q3", 38, 3.35),
    @(3, "2021-Q2", 33, 4.49),
    @(4, "2021-Q1", 4,  0.02)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = 2 + $i
    $totalRow = $totalRows[$i]
    $totalSheet.Cells.Item($r, 1).Value = $totalRow[0]
    $totalSheet.Cells.Item($r, 2).Value = $totalRow[1]
    $totalSheet.Cells.Item($r, 3).Value = $totalRow[2]
    $totalSheet.Cells.Item($r, 4).Value = $totalRow[3]
}

$srcTotalIndexCell = $template.Range("A2")
$srcTotalIndexCell.Copy()
$totalSheet.Range("A2:A6").PasteSpecial(-4122)

# Match the original sheet's page margins (0.75in/1in/0.5in) instead of the
# engine's generic new-sheet defaults.
$totalSheet.PageSetup.LeftMargin = 54
$totalSheet.PageSetup.RightMargin = 54
$totalSheet.PageSetup.TopMargin = 72
$totalSheet.PageSetup.BottomMargin = 72
$totalSheet.PageSetup.HeaderMargin = 36
$totalSheet.PageSetup.FooterMargin = 36

# Restore the originally active sheet/tab (2021-Q1) since adding sheets
# shifts selection onto the newest one.
$wb.Worksheets.Item(1).Activate()
